$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1750
$ws.Range("J19").Value = 1750
$ws.Range("L19").Value = 1750
$ws.Range("N19").Value = -2100

$ws.Range("H129").Value = 920.4483
$ws.Range("I129").Value = 311
$ws.Range("J129").Value = 1079.4348
$ws.Range("K129").Value = 933
$ws.Range("L129").Value = 3238.3044
$ws.Range("M129").Value = 4067
$ws.Range("N129").Value = -13238.3044

$ws.Range("H132").Value = 893302.0600000001
$ws.Range("I132").Value = 2423.9512
$ws.Range("J132").Value = 3502302.2
$ws.Range("K132").Value = 7271.8536
$ws.Range("L132").Value = 10506906.6
$ws.Range("M132").Value = -4741.8536
$ws.Range("N132").Value = -10511966.6

$ws.Range("H137").Value = 1640735.5
$ws.Range("I137").Value = 2174871.2
$ws.Range("J137").Value = 2719.1333
$ws.Range("K137").Value = 6524613.600000001
$ws.Range("L137").Value = 8157.3999
$ws.Range("M137").Value = -6522063.600000001
$ws.Range("N137").Value = -13257.3999

$ws.Range("H138").Value = 2606507
$ws.Range("I138").Value = 2007.5518
$ws.Range("J138").Value = 4764521
$ws.Range("K138").Value = 6022.6554
$ws.Range("L138").Value = 14293563
$ws.Range("M138").Value = -882.6553999999996
$ws.Range("N138").Value = -14303843

$ws.Range("H141").Value = 3001.3125
$ws.Range("I141").Value = 1583.5946
$ws.Range("J141").Value = 7770
$ws.Range("K141").Value = 4750.783799999999
$ws.Range("L141").Value = 23310
$ws.Range("M141").Value = 429.2162000000008
$ws.Range("N141").Value = -33670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2242.83
$ws.Range("I32").Value = 1684.5309
$ws.Range("J32").Value = 4622.9473
$ws.Range("K32").Value = 1684.5309
$ws.Range("L32").Value = 4622.9473
$ws.Range("M32").Value = -1397.5309
$ws.Range("N32").Value = -5196.9473

$ws.Range("H45").Value = 1550
$ws.Range("I45").Value = 1212.5
$ws.Range("J45").Value = 2900
$ws.Range("K45").Value = 1212.5
$ws.Range("L45").Value = 2900
$ws.Range("M45").Value = -835.5
$ws.Range("N45").Value = -3654

$ws.Range("H74").Value = 9316297
$ws.Range("I74").Value = 12383950
$ws.Range("J74").Value = 113340.445
$ws.Range("K74").Value = 12383950
$ws.Range("L74").Value = 113340.445
$ws.Range("M74").Value = -12383076
$ws.Range("N74").Value = -115088.445

$ws.Range("H77").Value = 9316297
$ws.Range("I77").Value = 12383950
$ws.Range("J77").Value = 113340.445
$ws.Range("K77").Value = 61919750
$ws.Range("L77").Value = 566702.2250000001
$ws.Range("M77").Value = -61915382
$ws.Range("N77").Value = -575438.2250000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1750.2745
$ws.Range("I134").Value = 1147.0278
$ws.Range("J134").Value = 3198.0667
$ws.Range("K134").Value = 3441.0834
$ws.Range("L134").Value = 9594.2001
$ws.Range("M134").Value = -906.0834000000004
$ws.Range("N134").Value = -14664.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 26317922
$ws.Range("I58").Value = 40002188
$ws.Range("J58").Value = 2024.7693
$ws.Range("K58").Value = 40002188
$ws.Range("L58").Value = 2024.7693
$ws.Range("M58").Value = -40001985
$ws.Range("N58").Value = -2430.7693

$ws.Range("H99").Value = 10066.667
$ws.Range("I99").Value = 8000
$ws.Range("K99").Value = 8000
$ws.Range("M99").Value = -6502

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H122").Value = 1854.5883
$ws.Range("J122").Value = 2134.7778
$ws.Range("L122").Value = 6404.3334
$ws.Range("N122").Value = -11304.3334

$ws.Range("H126").Value = 10066.667
$ws.Range("I126").Value = 8000
$ws.Range("K126").Value = 24000
$ws.Range("M126").Value = -21530

$ws.Range("H132").Value = 20951.02
$ws.Range("I132").Value = 1475.7693
$ws.Range("J132").Value = 79376.766
$ws.Range("K132").Value = 4427.3079
$ws.Range("L132").Value = 238130.298
$ws.Range("M132").Value = -1897.3079
$ws.Range("N132").Value = -243190.298

$ws.Range("H136").Value = 26317922
$ws.Range("I136").Value = 40002188
$ws.Range("J136").Value = 2024.7693
$ws.Range("K136").Value = 120006564
$ws.Range("L136").Value = 6074.3079
$ws.Range("M136").Value = -120004014
$ws.Range("N136").Value = -11174.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 512.8889
$ws.Range("I98").Value = 125
$ws.Range("J98").Value = 623.7143
$ws.Range("K98").Value = 375
$ws.Range("L98").Value = 1871.1429
$ws.Range("M98").Value = 1123
$ws.Range("N98").Value = -4867.1429

$ws.Range("H108").Value = 3651.0667
$ws.Range("I108").Value = 521
$ws.Range("J108").Value = 6389.875
$ws.Range("K108").Value = 1563
$ws.Range("L108").Value = 19169.625
$ws.Range("M108").Value = 1317
$ws.Range("N108").Value = -24929.625

$ws.Range("H109").Value = 2968.9167
$ws.Range("J109").Value = 3722.2222
$ws.Range("L109").Value = 11166.6666
$ws.Range("N109").Value = -13246.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1468.4615
$ws.Range("I113").Value = 1198
$ws.Range("J113").Value = 1637.5
$ws.Range("K113").Value = 1198
$ws.Range("L113").Value = 1637.5
$ws.Range("M113").Value = 972
$ws.Range("N113").Value = -5977.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2116.5
$ws.Range("I61").Value = 1999.6666
$ws.Range("J61").Value = 2233.3333
$ws.Range("K61").Value = 1999.6666
$ws.Range("L61").Value = 2233.3333
$ws.Range("M61").Value = -1797.6666
$ws.Range("N61").Value = -2637.3333

$ws.Range("H113").Value = 2116.5
$ws.Range("I113").Value = 1999.6666
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 1999.6666
$ws.Range("L113").Value = 2233.3333
$ws.Range("M113").Value = 170.3334
$ws.Range("N113").Value = -6573.3333

$ws.Range("H122").Value = 3246.4
$ws.Range("I122").Value = 2875.2144
$ws.Range("K122").Value = 8625.643199999999
$ws.Range("M122").Value = -6175.643199999999

$ws.Range("H132").Value = 46165.348
$ws.Range("I132").Value = 2320.375
$ws.Range("J132").Value = 146382.42
$ws.Range("K132").Value = 6961.125
$ws.Range("L132").Value = 439147.26
$ws.Range("M132").Value = -4431.125
$ws.Range("N132").Value = -444207.26

$ws.Range("H136").Value = 45726.312
$ws.Range("I136").Value = 35447.55
$ws.Range("J136").Value = 64356.562
$ws.Range("K136").Value = 106342.65
$ws.Range("L136").Value = 193069.686
$ws.Range("M136").Value = -106342.65
$ws.Range("N136").Value = -198169.686

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2425.375
$ws.Range("I96").Value = 2833.3333
$ws.Range("J96").Value = 2180.6
$ws.Range("K96").Value = 2833.3333
$ws.Range("L96").Value = 2180.6
$ws.Range("M96").Value = -1460.3333
$ws.Range("N96").Value = -4926.6

$ws.Range("H122").Value = 1820.0952
$ws.Range("I122").Value = 1286.3462
$ws.Range("K122").Value = 3859.0386
$ws.Range("M122").Value = -1409.0386

$ws.Range("H132").Value = 61981.605
$ws.Range("I132").Value = 44374.13
$ws.Range("J132").Value = 102478.8
$ws.Range("K132").Value = 133122.39
$ws.Range("L132").Value = 307436.4
$ws.Range("M132").Value = -130592.39
$ws.Range("N132").Value = -312496.4
